$wb = $excel.ActiveWorkbook

# Add the new sheet for Jan22 games, placed after the existing sheets
$ws2 = $wb.Worksheets.Item("Jan21(65pct)")
$ws3 = $wb.Worksheets.Add()
$ws3.Name = "Jan22(47pct)"
$ws3.Move($null, $ws2)

# Header row
$ws3.Range("A1").Value = "Away"
$ws3.Range("B1").Value = "Home"
$ws3.Range("C1").Value = "Away Score"
$ws3.Range("D1").Value = "Home Score"
$ws3.Range("A1:D1").Font.Bold = $true

$data = @(
  @("Wichita State", "South Florida", 41, 54),
  @("Duke", "Pittsburgh", 79, 64),
  @("Wake Forest", "Virginia", 45, 68),
  @("Clemson", "Florida State", 68, 77),
  @("Notre Dame", "Georgia Tech", 61, 63),
  @("Texas Tech", "Kansas State", 45, 58),
  @("Villanova", "Butler", 80, 72),
  @("Minnesota", "Michigan", 57, 59),
  @("Indiana", "Northwestern", 66, 73),
  @("Buffalo", "Northern Illinois", 75, 77),
  @("Central Michigan", "Akron", 67, 70),
  @("Miami (OH)", "Ball State", 71, 65),
  @("Bowling Green State", "Eastern Michigan", 80, 67),
  @("Toledo", "Kent State", 85, 87),
  @("Western Michigan", "Ohio", 76, 81),
  @("Boise State", "Air Force", 60, 74),
  @("New Mexico", "Nevada-Las Vegas", 58, 74),
  @("San Diego State", "Fresno State", 62, 66),
  @("Saint Peter's", "Niagara", 74, 72),
  @("Mississippi State", "Kentucky", 55, 76),
  @("Auburn", "South Carolina", 77, 80),
  @("Mississippi", "Alabama", 53, 74),
  @("Texas A&M", "Florida", 72, 81)
)

$row = 2
foreach ($game in $data) {
  $ws3.Cells.Item($row, 1).Value = $game[0]
  $ws3.Cells.Item($row, 2).Value = $game[1]
  $ws3.Cells.Item($row, 3).Value = $game[2]
  $ws3.Cells.Item($row, 4).Value = $game[3]
  $row = $row + 1
}

$ws3.Columns.Item(1).AutoFit()
$ws3.Columns.Item(2).AutoFit()

$ws3.Range("B11").Select()

$wb.Worksheets.Item("Jan21(65pct)").Range("D1:A1").Select()

$wb.Worksheets.Item("Jan22(47pct)").Activate()
